$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Highest Momentum Stocks")
$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A2").Value = "BMY"
$ws.Range("A3").Value = "CME"
$ws.Range("A4").Value = "K"
$ws.Range("A5").Value = "ABBV"
$ws.Range("A6").Value = "MO"
$ws.Range("A7").Value = "WEC"
$ws.Range("A8").Value = "KR"
$ws.Range("A9").Value = "VRSN"
$ws.Range("A10").Value = "GILD"
$ws.Range("A11").Value = "CBOE"

$ws = $wb.Worksheets.Item("Lowest PE Ratio Stocks")
$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A2").Value = "SITC"
$ws.Range("A3").Value = "WU"
$ws.Range("A4").Value = "AES"
$ws.Range("A5").Value = "MTDR"
$ws.Range("A6").Value = "F"
$ws.Range("A7").Value = "CALM"
$ws.Range("A8").Value = "HPE"
$ws.Range("A9").Value = "DVN"
$ws.Range("A10").Value = "GM"
$ws.Range("A11").Value = "APA"

$ws = $wb.Worksheets.Item("Highest Dividend Yield Stocks")
$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A2").Value = "SITC"
$ws.Range("A3").Value = "GMRE"
$ws.Range("A4").Value = "WBA"
$ws.Range("A5").Value = "WU"
$ws.Range("A6").Value = "LADR"
$ws.Range("A7").Value = "GOOD"
$ws.Range("A8").Value = "F"
$ws.Range("A9").Value = "DOW"
$ws.Range("A10").Value = "LYB"
$ws.Range("A11").Value = "KRNY"

$ws = $wb.Worksheets.Item("Most Stable Stocks")
$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A2").Value = "AKRO"
$ws.Range("A3").Value = "CALM"
$ws.Range("A4").Value = "STI"
$ws.Range("A5").Value = "BIIB"
$ws.Range("A6").Value = "GIS"
$ws.Range("A7").Value = "CPB"
$ws.Range("A8").Value = "GILD"
$ws.Range("A9").Value = "CTRA"
$ws.Range("A10").Value = "REGN"
$ws.Range("A11").Value = "SJM"

$ws = $wb.Worksheets.Item("Highest Market Cap Stocks")
$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A2").Value = "AAPL"
$ws.Range("A3").Value = "MSFT"
$ws.Range("A4").Value = "NVDA"
$ws.Range("A5").Value = "AMZN"
$ws.Range("A6").Value = "GOOGL"
$ws.Range("A7").Value = "META"
$ws.Range("A8").Value = "AVGO"
$ws.Range("A9").Value = "TSLA"
$ws.Range("A10").Value = "WMT"
$ws.Range("A11").Value = "JPM"

$excel.CutCopyMode = $false
Write-Host "done"